$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.425012350082397
$ws.Range("B1").Value = 3.200057744979858
$ws.Range("C1").Value = 1.946813225746155
$ws.Range("D1").Value = 1.343174695968628
$ws.Range("E1").Value = 1.11332094669342
